# "Blood Glucose Readings Sample.xlsx" - Add files via upload
# Adds 5 new days of readings above the existing table (2025-01-01..01-05),
# 2 new days below (2025-01-14..01-15), fills in a couple of previously-blank
# readings, converts Date/Day/Finger columns to fill-down formulas, drops the
# "Mounjaro 12.5" note, renames a note and renames two A1c labels on the Key
# sheet.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Readings")
$key = $wb.Worksheets.Item("Key")

# --- 1. Insert 5 new rows above row 2; old rows 2-9 shift down to rows 7-14 ---
$ws.Rows("2:6").Insert()

# New rows inherit the header's format by default - repaint them with the
# same look as the rest of the data rows (now row 7, the old row 2).
$ws.Range("A7:L7").Copy()
$ws.Range("A2:L6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Literal data for the 5 newly inserted rows (2025-01-01..01-05) ---
$ws.Range("A2").Value = 45658
$ws.Range("B2").Value = "Wednesday"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.3034722222222222
$ws.Range("E2").Value = 82
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 0.79305555555555551
$ws.Range("H2").Value = 123
$ws.Range("I2").Value = 2

$ws.Range("A3").Value = 45659
$ws.Range("D3").Value = 0.25833333333333336
$ws.Range("E3").Value = 90
$ws.Range("G3").Value = 0.7944444444444444
$ws.Range("H3").Value = 95

$ws.Range("A4").Value = 45660
$ws.Range("D4").Value = 0.26180555555555557
$ws.Range("E4").Value = 92
$ws.Range("G4").Value = 0.79791666666666672
$ws.Range("H4").Value = 89

$ws.Range("A5").Value = 45661
$ws.Range("D5").Value = 0.3034722222222222
$ws.Range("E5").Value = 86
$ws.Range("G5").Value = 0.79791666666666672
$ws.Range("H5").Value = 85

$ws.Range("A6").Value = 45662
$ws.Range("D6").Value = 0.35625000000000001
$ws.Range("E6").Value = 92
$ws.Range("G6").Value = 0.79513888888888884
$ws.Range("H6").Value = 121

# --- 3. Two new rows appended at the bottom (2025-01-14..01-15), cols A:I only ---
$ws.Range("A13:I13").Copy()
$ws.Range("A15:I16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A15").Value = 45671
$ws.Range("D15").Value = 0.25972222222222224
$ws.Range("E15").Value = 92
$ws.Range("G15").Value = 0.7944444444444444
$ws.Range("H15").Value = 89

$ws.Range("A16").Value = 45672

# --- 4. Fill in readings that were previously blank on what is now row 14 ---
$ws.Range("D13:E13").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$ws.Range("G13:H13").Copy()
$ws.Range("G14:H14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D14").Value = 0.26458333333333334
$ws.Range("E14").Value = 89
$ws.Range("G14").Value = 0.79722222222222228
$ws.Range("H14").Value = 114

# --- 5. Formulas: Date fills down, Day/Finger-rotation recompute from above ---
$ws.Range("A3:A16").Formula = "=A2+1"

for ($r = 3; $r -le 16; $r++) {
    $ws.Range("B$r").Formula = '=TEXT(A' + $r + ',"dddd")'
}

$ws.Range("C3:C16").Formula = "=IF(C2+1>4,1,C2+1)"
$ws.Range("F3:F16").Formula = "=IF(F2+1>4,1,F2+1)"
$ws.Range("I3:I16").Formula = "=IF(I2+1>4,1,I2+1)"

# --- 6. Drop the "Mounjaro 12.5" note that shifted onto row 9 ---
$ws.Range("L9").ClearContents()

# --- 7. Rename the "Extra reading" note (now on rows 11 and 12) ---
$ws.Range("L11").Value = "Extra reading to validate CGM reading"
$ws.Range("L12").Value = "Extra reading to validate CGM reading"

# --- 8. Rename the two A1c section headers on the Key sheet ---
$key.Range("A10").Value = "A1c (No Diabetes Present)"
$key.Range("E10").Value = "A1c (Diabetes Present)"

# --- 9. Tidy up selection to match the saved file ---
$ws.Range("A2").Select()

Write-Host "done"
